# Updates "想去人数" (column F) and, where applicable, "最低票价" (column G)
# values across the four sheets of the 杭州-漫展信息 workbook to reflect the
# newly scraped data (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

function Set-FG($Sheet, $Row, $F, $G) {
    $Sheet.Cells.Item($Row, 6).Value = $F
    if ($G -ne $null) {
        $Sheet.Cells.Item($Row, 7).Value = $G
    }
}

# ---------------------------------------------------------------------------
# Sheet "展览"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

Set-FG $ws1 3  2776
Set-FG $ws1 5  949
Set-FG $ws1 6  38
Set-FG $ws1 7  2950
Set-FG $ws1 8  1886
Set-FG $ws1 9  232
Set-FG $ws1 10 68
Set-FG $ws1 11 2541
Set-FG $ws1 12 573
Set-FG $ws1 13 261
Set-FG $ws1 14 5
Set-FG $ws1 16 137
Set-FG $ws1 17 125
Set-FG $ws1 18 9475
Set-FG $ws1 21 7439
Set-FG $ws1 22 11979
Set-FG $ws1 25 244
Set-FG $ws1 27 577
Set-FG $ws1 28 2695
Set-FG $ws1 29 240
Set-FG $ws1 31 2677
Set-FG $ws1 32 976
Set-FG $ws1 33 4
Set-FG $ws1 37 1069
Set-FG $ws1 38 26
Set-FG $ws1 39 362
Set-FG $ws1 40 59
Set-FG $ws1 41 562

# ---------------------------------------------------------------------------
# Sheet "演出"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

Set-FG $ws2 16 25 380
Set-FG $ws2 24 15

# ---------------------------------------------------------------------------
# Sheet "本地生活"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

Set-FG $ws3 2 634
Set-FG $ws3 4 184

# ---------------------------------------------------------------------------
# Sheet "全部类型"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

Set-FG $ws4 2  634
Set-FG $ws4 5  2776
Set-FG $ws4 8  949
Set-FG $ws4 9  38
Set-FG $ws4 11 2950
Set-FG $ws4 13 1886
Set-FG $ws4 14 232
Set-FG $ws4 15 2541
Set-FG $ws4 17 573
Set-FG $ws4 18 261
Set-FG $ws4 19 5
Set-FG $ws4 20 137
Set-FG $ws4 21 125
Set-FG $ws4 22 9475
Set-FG $ws4 25 7440
Set-FG $ws4 26 11979
Set-FG $ws4 29 244
Set-FG $ws4 32 577
Set-FG $ws4 34 2695
Set-FG $ws4 35 25 380
Set-FG $ws4 36 240
Set-FG $ws4 45 562
Set-FG $ws4 46 15
